# Apply the "To Do List" revision:
#  1. Delete the whole "为MiniMap添加边框（低）" task paragraph.
#  2. Move the hidden "_GoBack" bookmark from its old position
#     (inside the "管理、复现关键帧（高）" paragraph, right after
#     the "、") to just before the last remaining paragraph's run
#     ("雷雨等天气效果（低）") - i.e. where it ends up once the
#     "为MiniMap" paragraph and the now-empty numbering paragraph
#     collapse together.

$d = $word.ActiveDocument

# --- 1. Relocate the hidden _GoBack bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Delete the "为MiniMap添加边框（低）" paragraph entirely -------------
$rng = $d.Content
$found = $rng.Find.Execute("为MiniMap添加边框（低）", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $para.Range.Delete()
}

# --- 3. Re-create the bookmark right before the "雷雨等天气效果（低）" run --
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("雷雨等天气效果（低）", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $rng2) | Out-Null
}
